# "added new shop topo centras"
#
# Adds a new price-comparison row (topocentras.lt) to the "Iphone 11"
# table, and a corresponding "Error" row to the "usb atmintine 32gb"
# table (its scraper errored out for this run), growing both Excel
# Tables from A1:B3 to A1:B4.

$wb = $excel.ActiveWorkbook

# --- "Iphone 11" sheet / Table2: new Topo Centras offer -------------------
$wsIphone = $wb.Worksheets.Item("Iphone 11")
$tblIphone = $wsIphone.ListObjects.Item(1)
$tblIphone.ListRows.Add() | Out-Null

$wsIphone.Range("A4").Value = "529,,00"
$wsIphone.Range("B4").Value = "https://www.topocentras.lt/mobilusis-telefonas-apple-iphone-11-64gb-black.html"

# --- "usb atmintine 32gb" sheet / Table1: failed scrape row ---------------
$wsUsb = $wb.Worksheets.Item("usb atmintine 32gb")
$tblUsb = $wsUsb.ListObjects.Item(1)
$tblUsb.ListRows.Add() | Out-Null

$wsUsb.Range("A4").Value = "Error"
$wsUsb.Range("B4").Value = "Error"
